# NIT-9016108406.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# The worker KATRINY BRIGITH BOHORQUEZ HERNANDEZ (CC 1050960453) previously had
# a single overdue period (2508 / "2508") on the statement. A new overdue
# period (2509) is added for the same worker with the same amounts, which:
#   - inserts a new detail row (row 17) right below the existing one (row 16),
#     copying its formatting, and sets the period value to 2509
#   - doubles the total "VALOR MORA" (E11): 56940 -> 113880
#   - bumps "Cant. Periodos" (F13) from 1 to 2
# The two signature rows below the table are pushed down automatically by the
# row insert (old row 21/22 -> new row 22/23), matching the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old "NOMBRE DEL REPRESENTANTE LEGAL" block, right
# after the existing worker/period detail row, and clone that row's
# formatting into it.
$ws.Rows("17:17").Insert()
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# New period for the same worker (same CC/name/salary/value-in-arrears).
$ws.Range("E17").Value = "2509"

# Update the summary figures to reflect the second overdue period.
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2
